$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.08521032333374
$ws.Range("B1").Value = 2.625478982925415
$ws.Range("C1").Value = 1.930364012718201
$ws.Range("D1").Value = 1.788726091384888
$ws.Range("E1").Value = 1.664858341217041
